# aact-867: Update CDEK documentation based on Rebekah's input
$wb = $excel.ActiveWorkbook

$dataDefs = $wb.Worksheets.Item("Data Definitions")
$generalInfo = $wb.Worksheets.Item("General Info")

# Update the long description text on the "General Info" sheet (cell A2)
$newText = "`nClinical Drug Experience Knowledgebase (CDEK) is a database and web-platform to enable researchers and analysts to study all active pharmaceutical ingredients with evidence of clinical experience in humans. CDEK contains over 20,000 active pharmaceutical ingredients with surrounding annotated metadata (i.e. sponsoring organizations, indications, clinical trial data, chemical structure, etc). The Washington University team curated CDEK in part by disambiguating intervention and organization names from ClinicalTrials.gov (through AACT). This information was cross-referenced against entries in prominent international drug databases (e.g. PubChem, ChEMBL, DrugBank). The prototype web-application is designed to allow researchers of all backgrounds to access and query the data. It can be found at cdek.wustl.edu.`nOrganizations are often identified by a variety of names in ClinicalTrials.gov, making it difficult to correctly identify all studies associated with an organization. For example, ‘Washington University-St. Louis' and ‘Washington University School of Medicine' are both considered “Washington University in St. Louis”. In fact, Washington University in St. Louis and its schools/institutes, are represented 63 different ways in ClinicalTrials.gov.`nThe CDEK curation process currently focuses on organizations/affiliations identified as the sponsor, overall official or responsible parties for drug intervention trials in ClinicalTrials.gov as of September, 2017. (Organizations in the 'facilities' table were not included.) Each organization name was assigned a standardized 'preferred' (or parent) name, effectively creating an index of names that map the variety of names to a common name, making it easier to group studies by organization. `n"

$generalInfo.Range("A2").Value = $newText

# Selection / active-tab changes: the active sheet becomes "Data Definitions"
# (tab 0) rather than "General Info" (tab 1); each sheet's stored selection is
# updated to reflect the new state.
$dataDefs.Range("A2").Select()
$generalInfo.Range("A2").Select()
$dataDefs.Activate()

$wb.Save()
